# edit.ps1 - applies the "Recomendaciones" portfolio content swap described by the diff.
$d = $word.ActiveDocument

function Replace-ParagraphByFind {
    param([string]$FindText, [string]$NewParagraphXml)
    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph containing: $FindText"
    }
    $para = $rng.Paragraphs(1).Range
    $insertAt = $para.Duplicate
    $insertAt.End = $insertAt.Start
    $insertAt.InsertXML($NewParagraphXml)
}

function Delete-ParagraphByFind {
    param([string]$FindText)
    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph containing: $FindText"
    }
    $para = $rng.Paragraphs(1).Range
    $para.Delete()
}

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 6..9) Remove now-duplicated / superseded bullets FIRST (while their text is still
#           unique in the document, before later steps introduce duplicate wording) ---
Delete-ParagraphByFind "Desconfíe de aquellos e-mails"
Delete-ParagraphByFind "No propague aquellos mensajes de correo"

# --- 1) "Actualice regularmente ..." -> 3-run paragraph w/ proofErr around "ecommerce" ---
$p1Xml = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-PA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-PA"/></w:rPr><w:t xml:space="preserve">Con el creciente mundo tecnológico y los cambios en los hábitos de compra que ha impuesto la pandemia y la comodidad de estar en casa, muchos establecimientos comerciales se vieron obligados a acelerar estrategias para trasladar sus ventas al mundo de la internet. Las ganas de generar ventas llevaron a muchas empresas a implementar algunas estrategias de comercio electrónico o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-PA"/></w:rPr><w:t>ecommerce</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-PA"/></w:rPr><w:t>, sin considerar todo lo necesario para que este tipo de negocios sea totalmente “legal” y poder de esta forma satisfacer las necesidades y expectativas de los clientes.</w:t></w:r></w:p>
"@
Replace-ParagraphByFind "Actualice regularmente su sistema operativo" $p1Xml

# --- 2) "Instale un Antivirus ..." -> "Como, por ejemplo, las legislaciones ..." ---
$p2Xml = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-PA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-PA"/></w:rPr><w:t>Como, por ejemplo, las legislaciones que rigen este tipo de comercio, ya que muchas de las empresas y personas naturales y jurídicas desconocían incluso de la existencia de la dirección general de comercio electrónico como ente regulatorio de estas actividades dando paso a un sinfín de contratiempos que sumados a las fallas en los manejos logísticos de los inventarios llevaron el comercio electrónico en Panamá a la cuerda floja, sin embargo a partir de la pandemia se ha visto un gran crecimiento en la calidad y el atención ofrecida en la venta de productos y servicios regulados en internet, todo esto bajo las medidas tomadas por la DGCE y todas las facilidades que brinda en la actualidad para ejercer de esta manera</w:t></w:r></w:p>
"@
Replace-ParagraphByFind "Instale un Antivirus y actualícelo" $p2Xml

# --- 6) delete original "Ponga especial atención..." bullet BEFORE step 3 recreates its text ---
Delete-ParagraphByFind "Ponga especial atención en el tratamiento de su correo electrónico"

# --- 7) delete original "No abra mensajes..." bullet BEFORE step 4 recreates its text ---
Delete-ParagraphByFind "No abra mensajes de correo de remitentes desconocidos"

# --- 3) "Instale un Firewall ..." -> "Ponga especial atención en el tratamiento ..." ---
$p3Xml = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-PA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-PA"/></w:rPr><w:t>Ponga especial atención en el tratamiento de su correo electrónico, ya que es una de las herramientas más utilizadas para llevar a cabo estafas, introducir virus, etc.</w:t></w:r></w:p>
"@
Replace-ParagraphByFind "Instale un Firewall o Cortafuegos" $p3Xml

# --- 4) "tilice contraseñas seguras ... conveniente ... además ..." -> "No abra mensajes de correo de remitentes desconocidos." ---
$p4Xml = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-PA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-PA"/></w:rPr><w:t>No abra mensajes de correo de remitentes desconocidos.</w:t></w:r></w:p>
"@
Replace-ParagraphByFind "tilice contraseñas seguras" $p4Xml

# --- 5) "Navegue por páginas web seguras ..." -> "Mantenerse actualizado con el acontecer ciudadano ..." ---
$p5Xml = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="es-PA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-PA"/></w:rPr><w:t>Mantenerse actualizado con el acontecer ciudadano con medios como martesfinanciero.com</w:t></w:r></w:p>
"@
Replace-ParagraphByFind "Navegue por páginas web seguras" $p5Xml

Write-Host "All edits applied."
